$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (stored "width" = ColumnWidth + 5/6, so subtract 5/6 to land exactly) ---
$ws.Columns.Item(1).ColumnWidth = 11.1666667
$ws.Columns.Item(2).ColumnWidth = 13.1666667
$ws.Columns.Item(3).ColumnWidth = 13.1666667
$ws.Columns.Item(4).ColumnWidth = 23.1666667
$ws.Columns.Item(5).ColumnWidth = 24.1666667
$ws.Columns.Item(6).ColumnWidth = 24.1666667

# --- Row 1 header values ---
$ws.Range("A1").Value = "input_Name"
$ws.Range("B1").Value = "input_Name_1"
$ws.Range("C1").Value = "input_Name_2"
$ws.Range("D1").Value = "input_ReactSelectInput"
$ws.Range("E1").Value = "input_ReactSelectInput2"
$ws.Range("F1").Value = "input_ReactSelectInput3"

# Apply the same header style used by A1 to the new header cells
$ws.Range("A1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)

# --- Row 2: create empty (but present) text cells to mirror A2 ---
$ws.Range("B2").Value = "'"
$ws.Range("C2").Value = "'"
$ws.Range("D2").Value = "'"
$ws.Range("E2").Value = "'"
$ws.Range("F2").Value = "'"
$ws.Range("B2:F2").ClearFormats()
